$d = $word.ActiveDocument
$d.Content.Find.Execute("vo: Modifica di un Prodo", $true, $false, $false, $false, $false,
                         $true, 1, $false, "vo: aggiunta di un Prodo", 2)
